$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-30 Friday" "2026-01-31 Saturday"

Replace-Text "360÷5=" "589÷2="
Replace-Text "995÷8=" "568÷8="
Replace-Text "936÷6=" "412÷7="
Replace-Text "807÷9=" "792÷4="
Replace-Text "864÷6=" "711÷2="

Replace-Text "352÷6=" "281÷3="
Replace-Text "949÷4=" "736÷9="
Replace-Text "643÷6=" "105÷4="
Replace-Text "158÷2=" "118÷3="
Replace-Text "265÷4=" "507÷9="

Replace-Text "511÷6=" "698÷8="
Replace-Text "424÷6=" "827÷6="
Replace-Text "538÷8=" "444÷4="
Replace-Text "523÷5=" "238÷7="
Replace-Text "791÷2=" "221÷6="

Replace-Text "659÷6=" "904÷6="
Replace-Text "366÷7=" "908÷5="
Replace-Text "310÷4=" "467÷5="
Replace-Text "840÷7=" "277÷9="
Replace-Text "636÷6=" "774÷5="

Replace-Text "726÷2=" "645÷7="
Replace-Text "121÷3=" "842÷8="
Replace-Text "113÷3=" "526÷8="
Replace-Text "163÷4=" "580÷4="
Replace-Text "449÷3=" "440÷4="
